# Insert a new data row at sheet row 185 (pushing existing rows 185-216
# down to 186-217) and populate it with the new record, matching the
# commit's weekly price-update pattern for "Albahaca" at Vega Modelo de
# Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 185; this shifts every row
# from 185 onward down by one (old row 216 becomes row 217, etc.) and
# keeps their data untouched.
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new record.
$ws.Range("A185").Value = 10
$ws.Range("B185").Value = "Vega Modelo de Temuco"
$ws.Range("C185").Value = "La Araucanía"
$ws.Range("D185").Value = 44637
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = 100112052
$ws.Range("G185").Value = "Albahaca"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 80
$ws.Range("K185").Value = 7000
$ws.Range("L185").Value = 7000
$ws.Range("M185").Value = 7000
$ws.Range("N185").Value = "$/paquete"
$ws.Range("O185").Value = "Región de La Araucanía"
$ws.Range("P185").Value = 7000
$ws.Range("Q185").Value = 1
$ws.Range("R185").Value = "Hortaliza"
